$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '47.269.24'
$ws.Range("E2").Value = '  +2.25%  '
$ws.Range("D3").Value = '2.500.68'
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '323.52'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.80%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '108.97'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.31%  '
$ws.Range("E7").Value = '  +1.19%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.537'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("E10").Value = '  +8.73%  '
$ws.Range("E11").Value = '  +0.19%  '
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.39'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.20'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.43%  '
$ws.Range("D15").Value = '2.892.39'
$ws.Range("E15").Value = '  +1.65%  '
$ws.Range("D16").Value = '2.506.22'
$ws.Range("E16").Value = '  +1.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.854'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '47.214.94'
$ws.Range("E18").Value = '  +2.36%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.84'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("E20").Value = '  +3.25%  '
$ws.Range("E21").Value = '  +0.40%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.72'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +13.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.38%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '247.74'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("E25").Value = '  +3.15%  '
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("E28").Value = '  +3.86%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.07'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.65%  '
$ws.Range("E30").Value = '  +1.37%  '
$ws.Range("E31").Value = '  +7.53%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.84'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("E34").Value = '  +1.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0789'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.13%  '
$ws.Range("E36").Value = '  +0.23%  '
$ws.Range("E37").Value = '  +4.18%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.70'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '122.09'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.54%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '21.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '
$ws.Range("E44").Value = '  +2.09%  '
$ws.Range("D45").Value = '1.991.97'
$ws.Range("E45").Value = '  +0.92%  '
$ws.Range("E46").Value = '  +2.78%  '
$ws.Range("E47").Value = '  -1.52%  '
$ws.Range("E48").Value = '  -3.86%  '
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("E50").Value = '  +2.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.99'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.63%  '
